# Auto-generated bulk value update for Asura_Profits market-data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 16999.2
$ws.Range("J13").Value = 19996
$ws.Range("L13").Value = 19996
$ws.Range("N13").Value = -20334
$ws.Range("H15").Value = 798.35
$ws.Range("I15").Value = 798.35
$ws.Range("K15").Value = 2395.05
$ws.Range("M15").Value = -2226.05
$ws.Range("H41").Value = 302.05
$ws.Range("I41").Value = 350
$ws.Range("J41").Value = 270.08334
$ws.Range("K41").Value = 350
$ws.Range("L41").Value = 270.08334
$ws.Range("M41").Value = 90
$ws.Range("N41").Value = -1150.08334
$ws.Range("H93").Value = 44484
$ws.Range("J93").Value = 44484
$ws.Range("L93").Value = 44484
$ws.Range("N93").Value = -49476
$ws.Range("H112").Value = 1855.75
$ws.Range("J112").Value = 2048.8667
$ws.Range("L112").Value = 6146.6001
$ws.Range("N112").Value = -8362.6001
$ws.Range("H135").Value = 1358.1613
$ws.Range("I135").Value = 884.4231
$ws.Range("K135").Value = 7959.8079
$ws.Range("M135").Value = -5424.8079
$ws.Range("H141").Value = 14498.929
$ws.Range("I141").Value = 7921.923
$ws.Range("J141").Value = 100000
$ws.Range("K141").Value = 23765.769
$ws.Range("L141").Value = 300000
$ws.Range("M141").Value = -18585.769
$ws.Range("N141").Value = -310360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1280.2084
$ws.Range("I2").Value = 851.125
$ws.Range("K2").Value = 851.125
$ws.Range("M2").Value = -738.125
$ws.Range("H5").Value = 200
$ws.Range("I5").Value = 200
$ws.Range("K5").Value = 200
$ws.Range("M5").Value = -88
$ws.Range("H32").Value = 11649.304
$ws.Range("I32").Value = 12113.261
$ws.Range("J32").Value = 8448
$ws.Range("K32").Value = 12113.261
$ws.Range("L32").Value = 8448
$ws.Range("M32").Value = -11826.261
$ws.Range("N32").Value = -9022
$ws.Range("H61").Value = 1166.907
$ws.Range("I61").Value = 1113.5834
$ws.Range("J61").Value = 1234.2632
$ws.Range("K61").Value = 1113.5834
$ws.Range("L61").Value = 1234.2632
$ws.Range("M61").Value = -901.5834
$ws.Range("N61").Value = -1658.2632
$ws.Range("H110").Value = 2040.0769
$ws.Range("I110").Value = 1940.125
$ws.Range("J110").Value = 2200
$ws.Range("K110").Value = 1940.125
$ws.Range("L110").Value = 2200
$ws.Range("M110").Value = 104.875
$ws.Range("N110").Value = -6290
$ws.Range("H116").Value = 1280.2084
$ws.Range("I116").Value = 851.125
$ws.Range("K116").Value = 851.125
$ws.Range("M116").Value = 1442.875
$ws.Range("H132").Value = 1005309.7
$ws.Range("I132").Value = 1250877.2
$ws.Range("J132").Value = 23039.5
$ws.Range("K132").Value = 3752631.6
$ws.Range("L132").Value = 69118.5
$ws.Range("M132").Value = -3750101.6
$ws.Range("N132").Value = -74178.5
$ws.Range("H136").Value = 1166.907
$ws.Range("I136").Value = 1113.5834
$ws.Range("J136").Value = 1234.2632
$ws.Range("K136").Value = 3340.7502
$ws.Range("L136").Value = 3702.7896
$ws.Range("M136").Value = -790.7501999999999
$ws.Range("N136").Value = -8802.7896

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1280.2084
$ws.Range("I3").Value = 851.125
$ws.Range("K3").Value = 851.125
$ws.Range("M3").Value = -737.125
$ws.Range("H4").Value = 200
$ws.Range("I4").Value = 200
$ws.Range("K4").Value = 200
$ws.Range("M4").Value = -85
$ws.Range("H10").Value = 20250
$ws.Range("I10").Value = 20333.334
$ws.Range("J10").Value = 20000
$ws.Range("K10").Value = 20333.334
$ws.Range("L10").Value = 20000
$ws.Range("M10").Value = -20193.334
$ws.Range("N10").Value = -20280
$ws.Range("H109").Value = 23816.545
$ws.Range("J109").Value = 23816.545
$ws.Range("L109").Value = 23816.545
$ws.Range("N109").Value = -26590.545

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 625.7727
$ws.Range("I107").Value = 602
$ws.Range("J107").Value = 631.05554
$ws.Range("K107").Value = 602
$ws.Range("L107").Value = 631.05554
$ws.Range("M107").Value = 1318
$ws.Range("N107").Value = -4471.05554
$ws.Range("H122").Value = 1133.0952
$ws.Range("I122").Value = 1172.4667
$ws.Range("K122").Value = 3517.4001
$ws.Range("M122").Value = -1067.4001
$ws.Range("H132").Value = 2372.575
$ws.Range("I132").Value = 2196.8857
$ws.Range("J132").Value = 3602.4
$ws.Range("K132").Value = 6590.657099999999
$ws.Range("L132").Value = 10807.2
$ws.Range("M132").Value = -4060.657099999999
$ws.Range("N132").Value = -15867.2
$ws.Range("H134").Value = 1447.3243
$ws.Range("I134").Value = 1249.0714
$ws.Range("J134").Value = 2064.111
$ws.Range("K134").Value = 3747.2142
$ws.Range("L134").Value = 6192.333
$ws.Range("M134").Value = -1212.2142
$ws.Range("N134").Value = -11262.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3697.6155
$ws.Range("I39").Value = 700
$ws.Range("J39").Value = 4242.636
$ws.Range("K39").Value = 2100
$ws.Range("L39").Value = 12727.908
$ws.Range("M39").Value = -1806
$ws.Range("N39").Value = -13315.908
$ws.Range("H55").Value = 7779.4
$ws.Range("J55").Value = 8224.25
$ws.Range("L55").Value = 24672.75
$ws.Range("N55").Value = -25026.75
$ws.Range("H87").Value = 5981.6665
$ws.Range("I87").Value = 3972.5
$ws.Range("K87").Value = 11917.5
$ws.Range("M87").Value = -10669.5
$ws.Range("H90").Value = 5981.6665
$ws.Range("I90").Value = 3972.5
$ws.Range("K90").Value = 35752.5
$ws.Range("M90").Value = -29512.5
$ws.Range("H110").Value = 7660
$ws.Range("I110").Value = 5825
$ws.Range("K110").Value = 17475
$ws.Range("M110").Value = -13385
$ws.Range("H112").Value = 2494.7856
$ws.Range("I112").Value = 939.625
$ws.Range("J112").Value = 4568.3335
$ws.Range("K112").Value = 2818.875
$ws.Range("L112").Value = 13705.0005
$ws.Range("M112").Value = -1710.875
$ws.Range("N112").Value = -15921.0005
$ws.Range("H122").Value = 848.375
$ws.Range("I122").Value = 500.27274
$ws.Range("J122").Value = 1142.9231
$ws.Range("K122").Value = 4502.45466
$ws.Range("L122").Value = 10286.3079
$ws.Range("M122").Value = -2052.45466
$ws.Range("N122").Value = -15186.3079
$ws.Range("H123").Value = 1065
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H124").Value = 6705
$ws.Range("I124").Value = 1307.5
$ws.Range("K124").Value = 3922.5
$ws.Range("M124").Value = 987.5
$ws.Range("H129").Value = 1853391.9
$ws.Range("I129").Value = 743.3333
$ws.Range("J129").Value = 2084972.9
$ws.Range("K129").Value = 2229.9999
$ws.Range("L129").Value = 6254918.699999999
$ws.Range("M129").Value = 2770.0001
$ws.Range("N129").Value = -6264918.699999999
$ws.Range("H132").Value = 1817.2778
$ws.Range("I132").Value = 1166.7778
$ws.Range("K132").Value = 10501.0002
$ws.Range("M132").Value = -7971.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 21035
$ws.Range("J109").Value = 21035
$ws.Range("L109").Value = 21035
$ws.Range("N109").Value = -23115
$ws.Range("H122").Value = 4124.25
$ws.Range("J122").Value = 3061.6667
$ws.Range("L122").Value = 9185.000100000001
$ws.Range("N122").Value = -14085.0001
$ws.Range("H123").Value = 18889.533
$ws.Range("J123").Value = 18889.533
$ws.Range("L123").Value = 18889.533
$ws.Range("N123").Value = -23789.533

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 200
$ws.Range("I35").Value = 200
$ws.Range("K35").Value = 200
$ws.Range("M35").Value = 136
$ws.Range("H122").Value = 11369306
$ws.Range("I122").Value = 14711819
$ws.Range("J122").Value = 4761
$ws.Range("K122").Value = 44135457
$ws.Range("L122").Value = 14283
$ws.Range("M122").Value = -44133007
$ws.Range("N122").Value = -19183
$ws.Range("H131").Value = 21862.5
$ws.Range("J131").Value = 23369.889
$ws.Range("L131").Value = 23369.889
$ws.Range("N131").Value = -33449.889
$ws.Range("H136").Value = 3620.1458
$ws.Range("I136").Value = 3737.0715
$ws.Range("J136").Value = 2801.6667
$ws.Range("K136").Value = 11211.2145
$ws.Range("L136").Value = 8405.000100000001
$ws.Range("M136").Value = -8661.2145
$ws.Range("N136").Value = -13505.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1133.3334
$ws.Range("I113").Value = 934.7646999999999
$ws.Range("J113").Value = 1615.5714
$ws.Range("K113").Value = 2804.2941
$ws.Range("L113").Value = 4846.7142
$ws.Range("M113").Value = -634.2941000000001
$ws.Range("N113").Value = -9186.7142
$ws.Range("H126").Value = 3671.7693
$ws.Range("I126").Value = 4273.143
$ws.Range("J126").Value = 1146
$ws.Range("K126").Value = 12819.429
$ws.Range("L126").Value = 3438
$ws.Range("M126").Value = -10349.429
$ws.Range("N126").Value = -8378
$ws.Range("H132").Value = 2614.238
$ws.Range("I132").Value = 1825.9333
$ws.Range("K132").Value = 5477.7999
$ws.Range("M132").Value = -2947.7999
$ws.Range("H136").Value = 1265.9315
$ws.Range("I136").Value = 1185.1428
$ws.Range("K136").Value = 3555.4284
$ws.Range("M136").Value = -1005.4284
